$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '71.170.60'
$ws.Cells.Item(2, 5).Value = '  +0.36%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.863.23'
$ws.Cells.Item(3, 5).Value = '  +1.13%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '691.97'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +3.32%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '173.35'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +2.13%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.862.84'
$ws.Cells.Item(7, 5).Value = '  +1.17%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.05%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.12%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.19%  '

# Row 11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '7.45'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +6.75%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.461'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.70%  '

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '36.72'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.65%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.510.48'
$ws.Cells.Item(15, 5).Value = '  +0.97%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.869.92'
$ws.Cells.Item(16, 5).Value = '  +1.25%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '71.216.48'
$ws.Cells.Item(17, 5).Value = '  +0.51%  '

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '17.87'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.47%  '

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '7.26'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.83%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.42%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '11.16'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -5.13%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '488.42'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.54%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '0.721'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.41%  '

# Row 24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '84.84'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.81%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000147'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.66%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '12.40'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +1.14%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +1.29%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.91%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '4.016.03'
$ws.Cells.Item(29, 5).Value = '  +1.06%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.01%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '3.11'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +7.88%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '7.64'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +2.38%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.03%  '

# Row 34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '29.83'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.07%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.93%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '9.29'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.84%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '3.813.58'
$ws.Cells.Item(37, 5).Value = '  +1.04%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.01%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +1.58%  '

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +13.07%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '3.45'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.25%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.74%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +4.54%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.02%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.05%  '

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '163.95'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +4.01%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.000308'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +6.72%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Arweave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '44.73'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.92%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'OKB'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '48.71'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.26%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.304'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.97%  '

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '8.71'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +2.14%  '
